$wb = $excel.ActiveWorkbook

# --- Rename sheets: drop trailing/embedded periods from the "VEG"/"NON VEG" labels ---
$renames = @{
  "VEG. SOUP"                = "VEG SOUP"
  "NON. VEG. SOUP"            = "NON VEG SOUP"
  "SNACKS NON VEG."           = "SNACKS NON VEG"
  "SNACKS VEG."                = "SNACKS VEG"
  "TANDOORI VEG."              = "TANDOORI VEG"
  "INDIAN CURRIES VEG."        = "INDIAN CURRIES VEG"
  "TANDOOR NON VEG."           = "TANDOOR NON VEG"
  "INDIAN CURRIES NON VEG."    = "INDIAN CURRIES NON VEG"
}
foreach ($oldName in $renames.Keys) {
  $ws = $wb.Worksheets.Item($oldName)
  $ws.Name = $renames[$oldName]
}

# --- Add a new menu item ("EGG Pakoda" / 70) as the next row on the SNACKS VEG sheet ---
$snacksVeg = $wb.Worksheets.Item("SNACKS VEG")

# Copy the formatting of the preceding row (44) down onto the new row (45) first ...
$snacksVeg.Range("A44:B44").Copy()
$snacksVeg.Range("A45:B45").PasteSpecial(-4122)

# ... then fill in the new values.
$snacksVeg.Range("A45").Value = "EGG Pakoda"
$snacksVeg.Range("B45").Value = 70

# --- Update cursor/selection state on a few sheets, matching the edited session ---
$vegSoup = $wb.Worksheets.Item("VEG SOUP")
$vegSoup.Activate()
$vegSoup.Range("C29").Select()

$snacksVeg.Activate()
$snacksVeg.Range("A46").Select()

$tandoorNonVeg = $wb.Worksheets.Item("TANDOOR NON VEG")
$tandoorNonVeg.Activate()
$tandoorNonVeg.Range("B17").Select()

$bread = $wb.Worksheets.Item("BREAD")
$bread.Activate()
$bread.Range("G8").Select()
